$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 191
$ws.Range("A191").Value = 190
$ws.Range("B191").Value = "Monday, Jan 16"
$ws.Range("C191").Value = "6:10 AM"
$ws.Range("D191").Value = "FR2670"
$ws.Range("E191").Value = "London"
$ws.Range("F191").Value = "(STN)"
$ws.Range("G191").Value = "Buzz "
$ws.Range("H191").Value = "B38M"
$ws.Range("I191").Value = "(SP-RZE)"
$ws.Range("J191").Value = "6:12 AM"
$ws.Range("L191").Value = "0 hours, 2 minutes"

# Row 192
$ws.Range("A192").Value = 191
$ws.Range("B192").Value = "Monday, Jan 16"
$ws.Range("C192").Value = "6:50 AM"
$ws.Range("D192").Value = "FR1963"
$ws.Range("E192").Value = "Milan"
$ws.Range("F192").Value = "(BGY)"
$ws.Range("G192").Value = "Ryanair "
$ws.Range("H192").Value = "B738"
$ws.Range("I192").Value = "(SP-RKD)"
$ws.Range("J192").Value = "7:00 AM"
$ws.Range("L192").Value = "0 hours, 10 minutes"

# Row 193
$ws.Range("A193").Value = 192
$ws.Range("B193").Value = "Monday, Jan 16"
$ws.Range("C193").Value = "7:00 AM"
$ws.Range("D193").Value = "FR1115"
$ws.Range("E193").Value = "Rome"
$ws.Range("F193").Value = "(CIA)"
$ws.Range("G193").Value = "Ryanair "
$ws.Range("H193").Value = "B738"
$ws.Range("I193").Value = "(SP-RKF)"
$ws.Range("J193").Value = "6:58 AM"
$ws.Range("L193").Value = "0 hours, -2 minutes"

# Row 194
$ws.Range("A194").Value = 193
$ws.Range("B194").Value = "Monday, Jan 16"
$ws.Range("C194").Value = "7:00 AM"
$ws.Range("D194").Value = "FR4533"
$ws.Range("E194").Value = "Porto"
$ws.Range("F194").Value = "(OPO)"
$ws.Range("G194").Value = "Ryanair "
$ws.Range("H194").Value = "B38M"
$ws.Range("I194").Value = "(SP-RZO)"
$ws.Range("J194").Value = "7:18 AM"
$ws.Range("L194").Value = "0 hours, 18 minutes"

# Row 195
$ws.Range("A195").Value = 194
$ws.Range("B195").Value = "Monday, Jan 16"
$ws.Range("C195").Value = "7:05 AM"
$ws.Range("D195").Value = "FR6121"
$ws.Range("E195").Value = "Tenerife"
$ws.Range("F195").Value = "(TFS)"
$ws.Range("G195").Value = "Ryanair "
$ws.Range("H195").Value = "B738"
$ws.Range("I195").Value = "(SP-RSP)"
$ws.Range("J195").Value = "7:33 AM"
$ws.Range("L195").Value = "0 hours, 28 minutes"

# Row 196
$ws.Range("A196").Value = 195
$ws.Range("B196").Value = "Monday, Jan 16"
$ws.Range("C196").Value = "7:30 AM"
$ws.Range("D196").Value = "UNKNOWN"
$ws.Range("E196").Value = "Nice"
$ws.Range("F196").Value = "(NCE)"
$ws.Range("G196").Value = "NetJets Europe "
$ws.Range("H196").Value = "F2TH"
$ws.Range("I196").Value = "(CS-DLF)"
$ws.Range("J196").Value = "8:02 AM"
$ws.Range("L196").Value = "0 hours, 32 minutes"

# Row 197
$ws.Range("A197").Value = 196
$ws.Range("B197").Value = "Monday, Jan 16"
$ws.Range("C197").Value = "7:45 AM"
$ws.Range("D197").Value = "FR2000"
$ws.Range("E197").Value = "Cologne"
$ws.Range("F197").Value = "(CGN)"
$ws.Range("G197").Value = "Buzz "
$ws.Range("H197").Value = "B38M"
$ws.Range("I197").Value = "(SP-RZG)"
$ws.Range("J197").Value = "7:59 AM"
$ws.Range("L197").Value = "0 hours, 14 minutes"

# Row 198
$ws.Range("A198").Value = 197
$ws.Range("B198").Value = "Monday, Jan 16"
$ws.Range("C198").Value = "8:50 AM"
$ws.Range("D198").Value = "9U424"
$ws.Range("E198").Value = "Chisinau"
$ws.Range("F198").Value = "(KIV)"
$ws.Range("G198").Value = "Air Moldova "
$ws.Range("H198").Value = "A319"
$ws.Range("I198").Value = "(ER-AXL)"
$ws.Range("J198").Value = "8:46 AM"
$ws.Range("L198").Value = "0 hours, -4 minutes"

# Row 199
$ws.Range("A199").Value = 198
$ws.Range("B199").Value = "Monday, Jan 16"
$ws.Range("C199").Value = "9:00 AM"
$ws.Range("D199").Value = "FR4525"
$ws.Range("E199").Value = "Edinburgh"
$ws.Range("F199").Value = "(EDI)"
$ws.Range("G199").Value = "Ryanair "
$ws.Range("H199").Value = "B738"
$ws.Range("I199").Value = "(SP-RKP)"
$ws.Range("J199").Value = "9:05 AM"
$ws.Range("L199").Value = "0 hours, 5 minutes"

# Row 200
$ws.Range("A200").Value = 199
$ws.Range("B200").Value = "Monday, Jan 16"
$ws.Range("C200").Value = "9:35 AM"
$ws.Range("D200").Value = "FR4670"
$ws.Range("E200").Value = "Stockholm"
$ws.Range("F200").Value = "(ARN)"
$ws.Range("G200").Value = "Ryanair "
$ws.Range("H200").Value = "B738"
$ws.Range("I200").Value = "(9H-QDK)"
$ws.Range("J200").Value = "9:48 AM"
$ws.Range("L200").Value = "0 hours, 13 minutes"
